# This script normalizes a handful of game-title entries in column A of
# Sheet1 by removing punctuation/suffixes ("Remake", colons, trailing
# spaces, etc.) so that duplicate/near-duplicate titles collapse to a
# single canonical spelling.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A13").Value = "Dragon Ball Z Kakarot"
$ws.Range("A10").Value = "Crash Bandicoot"
$ws.Range("A49").Value = "Resident Evil 4"
$ws.Range("A46").Value = "Resident Evil 2"
$ws.Range("A47").Value = "Resident Evil 2"
$ws.Range("A60").Value = "Spider-Man Miles Morales"
$ws.Range("A61").Value = "Spider-Man Miles Morales"
$ws.Range("A68").Value = "The Witcher 3 Wild Hunt"
$ws.Range("A69").Value = "The Witcher 3 Wild Hunt"
$ws.Range("A71").Value = "Uncharted 4 A Thief's End"
$ws.Range("A76").Value = "Mortal Kombat 1"

# Reflect the scroll position / active selection recorded for the sheet
# after the edit.
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("A77:B78").Select()
$ws.Range("B77").Activate()
